$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LP Rewards")

# Format the full data range as text up front so that numeric-looking
# strings (token addresses, decimal amounts, etc.) are preserved exactly
# as text, matching how the source data is stored.
$ws.Range("A1:E25").NumberFormat = "@"

$rows = @(
    ,@("0xc1612C97537c2CC62a11FC4516367AB6F62d4B23", "0.000001673806345634", "0.8", "6.89", "160")
    ,@("0xa4fF2Ff7c9dCc78b76Ac381E2d06938B88bd8804", "0.003594689476742665", "1738.04", "14969.18", "347608")
    ,@("0xCD7a60306d8315B0b611De60B38Ae0B623062b2e", "0.000303845582308413", "183.93", "1424.7", "33084")
    ,@("0x248784Af80094fe92F915674b2f4409366443562", "0.005060191022620786", "2189.57", "0", "0")
    ,@("0x93E7fB017883f56fBc43f5723eCcCf6bF90FE6cF", "0.000364553713849458", "140.84", "1365.54", "31710")
    ,@("0x8F1c51E98Af7C7dbB24654ACB05781E0e96e008F", "0.082142650014355746", "42455.63", "353860.12", "8217188")
    ,@("0x4651B2E464e94c730725F3D5B2737418E2e5fd32", "0.019648622723011631", "9500.17", "81821.94", "1900034")
    ,@("0x0380ad3322Df94334C2f30CEE24D3086FC6F3445", "0.00381512308505338", "1268.64", "0", "0")
    ,@("0x2fe6f7c52EccC8Fea03289EC10213033a7c364b7", "0.000276182821902365", "61.39", "0", "0")
    ,@("0xdcE6F2ccaf54de860a5f80DB7503154957d216D5", "0.011746297729462631", "5386.43", "0", "0")
    ,@("0xBa44825f3215A71b34a702A70373afF8d4A52390", "0.009905592850166421", "5029", "42281.29", "981838")
    ,@("0xa614A83132a2e7368aDa71dA9331817c33706770", "0", "0", "0", "0")
    ,@("0xc38eDb1B53c7c034867BdE801e24826687F24466", "0.008026972749126571", "4457.34", "7164.71", "166376")
    ,@("0xabF76381Ab5538598AA4cd77d01De16c3F48a54A", "0.003697156090803694", "1259.12", "0", "0")
    ,@("0xf51435DD16E9a346f016A0a9123A7A87D80b28AA", "0.003600962835922976", "1922.75", "15777.65", "366382")
    ,@("0xf5DFA78573d859Ba64522B44FF0F7f970cFc5D78", "0.011698625575774027", "4055.78", "0", "0")
    ,@("0x5196E6ff6b95536874ba6405c11f465e026099e0", "0.00118886197213139", "574.7", "4950.18", "114951")
    ,@("0x19D27Bd32c639C4b5DaA9b9D05AD2401A79224F0", "0", "0", "0", "0")
    ,@("0x87029960cb6d1011289d35bC8944d37B88C46B72", "0.001735955596785356", "764.19", "4973.43", "115491")
    ,@("0x3663b2EE275616bfA14e77c69dc02b79e6C580F0", "0.005049157278053823", "3010.46", "23477.05", "545174")
    ,@("0x92119de4787cb947A4A93FEEC7c19C3AFe1b5729", "0.123174058999187785", "60862.03", "0", "0")
    ,@("0xCDA0c8a5f22565ae0efB889E5dBb6CE60498d0b4", "0.01020271109981372", "4286.15", "0", "0")
    ,@("0xe3B1ebcD19c28459aD75B44C13095cA93E53cbB0", "0.00092088557386215", "387.47", "3585.97", "83272")
    ,@("0xf6786b4edFa78d8d5cbfB93A97fF260e873152Be", "0.023721435439256031", "10008.11", "92489.43", "2147750")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowNum = $i + 2
    $data = New-Object "object[,]" 1,5
    for ($j = 0; $j -lt 5; $j++) {
        $data[0,$j] = $rows[$i][$j]
    }
    $ws.Range($ws.Cells.Item($rowNum,1), $ws.Cells.Item($rowNum,5)).Value = $data
}
